$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Column P (ShipmentTracking) - new tracking numbers for every row
Set-TextValue "P2"  "320018655634"
Set-TextValue "P3"  "320018655645"
Set-TextValue "P4"  "320018655678"
Set-TextValue "P5"  "320018655690"
Set-TextValue "P6"  "320018655760"
Set-TextValue "P7"  "320018655781"
Set-TextValue "P8"  "320018655818"
Set-TextValue "P9"  "320018655830"
Set-TextValue "P10" "320018655862"
Set-TextValue "P11" "320018655884"
Set-TextValue "P12" "320018655921"
Set-TextValue "P13" "320018655943"
Set-TextValue "P14" "320018655976"
Set-TextValue "P15" "320018655998"
Set-TextValue "P16" "320018645839"
Set-TextValue "P17" "320018645850"
Set-TextValue "P18" "320018645894"
Set-TextValue "P19" "320018645910"
Set-TextValue "P20" "320018645942"
Set-TextValue "P21" "320018645964"
Set-TextValue "P22" "320018645997"
Set-TextValue "P23" "320018646000"
Set-TextValue "P24" "320018646011"
Set-TextValue "P25" "320018646022"
Set-TextValue "P26" "320018646033"

# Column Q (ExpectedRate) and R (Result) for the rows whose outcome changed
Set-TextValue "Q4"  "$74.03"
Set-TextValue "R4"  "FAIL"

Set-TextValue "Q22" "$223.37"
Set-TextValue "R22" "PASS"

Set-TextValue "Q23" "$436.98"
Set-TextValue "R23" "PASS"

Set-TextValue "Q24" "$248.51"

Set-TextValue "Q25" "$52.88"
Set-TextValue "R25" "PASS"

Set-TextValue "Q26" "$1,171.41"
Set-TextValue "R26" "PASS"
